$d = $word.ActiveDocument

$d.Content.Find.Execute("95÷4=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "77÷2=38, 1", 2) | Out-Null
$d.Content.Find.Execute("89÷6=14, 5", $true, $false, $false, $false, $false, $true, 1, $false, "91÷8=11, 3", 2) | Out-Null
$d.Content.Find.Execute("74÷3=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "65÷8=8, 1", 2) | Out-Null
$d.Content.Find.Execute("28÷4=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2) | Out-Null
$d.Content.Find.Execute("23÷5=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2) | Out-Null
$d.Content.Find.Execute("42÷6=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=4, 0", 2) | Out-Null
$d.Content.Find.Execute("20÷7=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "87÷3=29, 0", 2) | Out-Null
$d.Content.Find.Execute("93÷7=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=11, 3", 2) | Out-Null
$d.Content.Find.Execute("64÷7=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷5=7, 4", 2) | Out-Null
$d.Content.Find.Execute("55÷5=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷2=28, 0", 2) | Out-Null
$d.Content.Find.Execute("26÷6=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=19, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷4=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2) | Out-Null
$d.Content.Find.Execute("11÷8=1, 3", $true, $false, $false, $false, $false, $true, 1, $false, "89÷8=11, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=5, 3", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("62÷7=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "86÷4=21, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=25, 0", 2) | Out-Null
$d.Content.Find.Execute("24÷9=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "80÷9=8, 8", 2) | Out-Null
$d.Content.Find.Execute("61÷8=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "77÷2=38, 1", 2) | Out-Null
$d.Content.Find.Execute("68÷5=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=12, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷8=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=7, 7", 2) | Out-Null
$d.Content.Find.Execute("41÷8=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "68÷3=22, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷8=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=6, 0", 2) | Out-Null
$d.Content.Find.Execute("16÷9=1, 7", $true, $false, $false, $false, $false, $true, 1, $false, "34÷6=5, 4", 2) | Out-Null
$d.Content.Find.Execute("92÷4=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "21÷3=7, 0", 2) | Out-Null
